# Fix typo in KCDC mayor candidate's name: "Guruanthan K (Guru)" -> "Gurunathan K (Guru)"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldName = "Guruanthan K (Guru)"
$newName = "Gurunathan K (Guru)"

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value2 -eq $oldName) {
        $cell.Value2 = $newName
    }
}
